# "Generate Report for Handback"
#
# Row 7 of both locale sheets (zh-cn, de-de) previously had no target /
# handback information for the 2a54b322-... file. A new handback report
# run fills in:
#   - I7  Latest Target File      -> hyperlinked to the source .md (same
#                                     text/link as column A's hyperlink)
#   - J7  Latest Handback File    -> the generated handback xliff file name
#                                     (same value already shown in column G)
#   - K7  Latest Handback DateTime-> the new handback timestamp
#   - P7  Error Detail            -> a "handback file is not the latest"
#                                     warning, identical wording/urls on
#                                     both locale sheets
# plus a new external hyperlink on I7 pointing at the same "latest" commit
# URL already used by column A's link (rId11 / rId8).

$wb = $excel.ActiveWorkbook

$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/22738884b7832339f3578cb10d510bd7809fdfc6/e2e/2a54b322-b825-4b04-a1ae-21c21727229c.md"
$mdDisplay = "2a54b322-b825-4b04-a1ae-21c21727229c.md"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c40bde7ce4947dafdd154165c53e85358cc863fa/e2e/2a54b322-b825-4b04-a1ae-21c21727229c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/22738884b7832339f3578cb10d510bd7809fdfc6/e2e/2a54b322-b825-4b04-a1ae-21c21727229c.md."

# ---------------------------------------------------------------------
# zh-cn sheet, row 7
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I7").Value = $mdDisplay
$wsZh.Range("I7").Font.Underline = 2
$wsZh.Range("I7").Font.Color = 15570276
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $mdUrl, "", "", $mdDisplay)

$wsZh.Range("J7").Value = "2a54b322-b825-4b04-a1ae-21c21727229c.d2c2fcd6fd81c4338ad71302c20cc37dac9a4886.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-09-03 02:59:11"
$wsZh.Range("K7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("P7").Value = $errorDetail

# ---------------------------------------------------------------------
# de-de sheet, row 7
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I7").Value = $mdDisplay
$wsDe.Range("I7").Font.Underline = 2
$wsDe.Range("I7").Font.Color = 15570276
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $mdUrl, "", "", $mdDisplay)

$wsDe.Range("J7").Value = "2a54b322-b825-4b04-a1ae-21c21727229c.d2c2fcd6fd81c4338ad71302c20cc37dac9a4886.de-de.xlf"
$wsDe.Range("K7").Value = "2016-09-03 02:59:18"
$wsDe.Range("K7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("P7").Value = $errorDetail
